$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Part 1: fix row ordering for duplicate-timestamp match pairs ---
# (content in columns F:V swapped/rotated between rows sharing the same kickoff date,
#  columns A:E are unchanged)

# Row 8
$ws.Cells.Item(8, 6).Value = 'Elana Torun'
$ws.Cells.Item(8, 7).Value = 3
$ws.Cells.Item(8, 8).Value = 'Blekitni Stargard'
$ws.Cells.Item(8, 9).Value = 0
$ws.Cells.Item(8, 10).Value = 2.12
$ws.Cells.Item(8, 11).Value = '05/08/2023 16:50'
$ws.Cells.Item(8, 12).Value = 2.26
$ws.Cells.Item(8, 13).Value = '05/08/2023 16:59'
$ws.Cells.Item(8, 14).Value = 3.38
$ws.Cells.Item(8, 15).Value = '05/08/2023 16:50'
$ws.Cells.Item(8, 16).Value = 3.34
$ws.Cells.Item(8, 17).Value = '05/08/2023 16:59'
$ws.Cells.Item(8, 18).Value = 2.96
$ws.Cells.Item(8, 19).Value = '05/08/2023 16:50'
$ws.Cells.Item(8, 20).Value = 2.75
$ws.Cells.Item(8, 21).Value = '05/08/2023 16:59'
$ws.Cells.Item(8, 22).Value = 'https://www.betexplorer.com/football/poland/iii-liga-group-ii/elana-torun-blekitni-stargard/KxcO4DZR/'

# Row 9
$ws.Cells.Item(9, 6).Value = 'Swinoujscie'
$ws.Cells.Item(9, 7).Value = 0
$ws.Cells.Item(9, 8).Value = 'Nowe Skalmierzyce'
$ws.Cells.Item(9, 9).Value = 3
$ws.Cells.Item(9, 10).Value = 2.12
$ws.Cells.Item(9, 11).Value = '05/08/2023 16:50'
$ws.Cells.Item(9, 12).Value = 2.12
$ws.Cells.Item(9, 13).Value = '05/08/2023 16:50'
$ws.Cells.Item(9, 14).Value = 3.38
$ws.Cells.Item(9, 15).Value = '05/08/2023 16:50'
$ws.Cells.Item(9, 16).Value = 3.38
$ws.Cells.Item(9, 17).Value = '05/08/2023 16:50'
$ws.Cells.Item(9, 18).Value = 2.96
$ws.Cells.Item(9, 19).Value = '05/08/2023 16:50'
$ws.Cells.Item(9, 20).Value = 2.96
$ws.Cells.Item(9, 21).Value = '05/08/2023 16:50'
$ws.Cells.Item(9, 22).Value = 'https://www.betexplorer.com/football/poland/iii-liga-group-ii/swinoujscie-nowe-skalmierzyce/jLmfaick/'

# Row 13
$ws.Cells.Item(13, 6).Value = 'Blekitni Stargard'
$ws.Cells.Item(13, 7).Value = 3
$ws.Cells.Item(13, 8).Value = 'Solec Kujawski'
$ws.Cells.Item(13, 9).Value = 1
$ws.Cells.Item(13, 10).Value = 1.95
$ws.Cells.Item(13, 11).Value = '11/08/2023 03:13'
$ws.Cells.Item(13, 12).Value = 2.16
$ws.Cells.Item(13, 13).Value = '12/08/2023 14:59'
$ws.Cells.Item(13, 14).Value = 3.28
$ws.Cells.Item(13, 15).Value = '11/08/2023 03:13'
$ws.Cells.Item(13, 16).Value = 3.86
$ws.Cells.Item(13, 17).Value = '12/08/2023 14:59'
$ws.Cells.Item(13, 18).Value = 3.02
$ws.Cells.Item(13, 19).Value = '11/08/2023 03:13'
$ws.Cells.Item(13, 20).Value = 2.61
$ws.Cells.Item(13, 21).Value = '12/08/2023 14:59'
$ws.Cells.Item(13, 22).Value = 'https://www.betexplorer.com/football/poland/iii-liga-group-ii/blekitni-stargard-unia-solec-kujawski/h6eQ1pJm/'

# Row 14
$ws.Cells.Item(14, 6).Value = 'Kleczew'
$ws.Cells.Item(14, 7).Value = 3
$ws.Cells.Item(14, 8).Value = 'Pogon Szczecin II'
$ws.Cells.Item(14, 9).Value = 2
$ws.Cells.Item(14, 10).Value = 3.51
$ws.Cells.Item(14, 11).Value = '11/08/2023 03:13'
$ws.Cells.Item(14, 12).Value = 1.99
$ws.Cells.Item(14, 13).Value = '12/08/2023 14:58'
$ws.Cells.Item(14, 14).Value = 3.51
$ws.Cells.Item(14, 15).Value = '11/08/2023 03:13'
$ws.Cells.Item(14, 16).Value = 3.7
$ws.Cells.Item(14, 17).Value = '12/08/2023 14:58'
$ws.Cells.Item(14, 18).Value = 1.76
$ws.Cells.Item(14, 19).Value = '11/08/2023 03:13'
$ws.Cells.Item(14, 20).Value = 2.89
$ws.Cells.Item(14, 21).Value = '12/08/2023 14:58'
$ws.Cells.Item(14, 22).Value = 'https://www.betexplorer.com/football/poland/iii-liga-group-ii/kleczew-pogon-szczecin/QiupME42/'

# Row 20
$ws.Cells.Item(20, 6).Value = 'Solec Kujawski'
$ws.Cells.Item(20, 7).Value = 0
$ws.Cells.Item(20, 8).Value = 'Zawisza'
$ws.Cells.Item(20, 9).Value = 2
$ws.Cells.Item(20, 10).Value = 3.48
$ws.Cells.Item(20, 11).Value = '18/08/2023 02:14'
$ws.Cells.Item(20, 12).Value = 5.06
$ws.Cells.Item(20, 13).Value = '19/08/2023 13:59'
$ws.Cells.Item(20, 14).Value = 3.41
$ws.Cells.Item(20, 15).Value = '18/08/2023 02:14'
$ws.Cells.Item(20, 16).Value = 4.28
$ws.Cells.Item(20, 17).Value = '19/08/2023 13:59'
$ws.Cells.Item(20, 18).Value = 1.75
$ws.Cells.Item(20, 19).Value = '18/08/2023 02:14'
$ws.Cells.Item(20, 20).Value = 1.48
$ws.Cells.Item(20, 21).Value = '19/08/2023 13:59'
$ws.Cells.Item(20, 22).Value = 'https://www.betexplorer.com/football/poland/iii-liga-group-ii/unia-solec-kujawski-zawisza/p0K3F3mm/'

# Row 21
$ws.Cells.Item(21, 6).Value = 'Notec Czarnkow'
$ws.Cells.Item(21, 7).Value = 1
$ws.Cells.Item(21, 8).Value = 'Swit Skolwin'
$ws.Cells.Item(21, 9).Value = 5
$ws.Cells.Item(21, 10).Value = 2.86
$ws.Cells.Item(21, 11).Value = '19/08/2023 09:25'
$ws.Cells.Item(21, 12).Value = 3.05
$ws.Cells.Item(21, 13).Value = '19/08/2023 13:02'
$ws.Cells.Item(21, 14).Value = 3.43
$ws.Cells.Item(21, 15).Value = '19/08/2023 09:25'
$ws.Cells.Item(21, 16).Value = 3.48
$ws.Cells.Item(21, 17).Value = '19/08/2023 13:02'
$ws.Cells.Item(21, 18).Value = 2.13
$ws.Cells.Item(21, 19).Value = '19/08/2023 09:25'
$ws.Cells.Item(21, 20).Value = 2.04
$ws.Cells.Item(21, 21).Value = '19/08/2023 13:02'
$ws.Cells.Item(21, 22).Value = 'https://www.betexplorer.com/football/poland/iii-liga-group-ii/notec-czarnkow-swit-skolwin/zRUTL5AJ/'

# Row 24
$ws.Cells.Item(24, 6).Value = 'Swinoujscie'
$ws.Cells.Item(24, 7).Value = 0
$ws.Cells.Item(24, 8).Value = 'Blekitni Stargard'
$ws.Cells.Item(24, 9).Value = 3
$ws.Cells.Item(24, 10).Value = 2.32
$ws.Cells.Item(24, 11).Value = '19/08/2023 09:26'
$ws.Cells.Item(24, 12).Value = 2.34
$ws.Cells.Item(24, 13).Value = '19/08/2023 16:59'
$ws.Cells.Item(24, 14).Value = 3.3
$ws.Cells.Item(24, 15).Value = '19/08/2023 09:26'
$ws.Cells.Item(24, 16).Value = 3.33
$ws.Cells.Item(24, 17).Value = '19/08/2023 16:59'
$ws.Cells.Item(24, 18).Value = 2.64
$ws.Cells.Item(24, 19).Value = '19/08/2023 09:26'
$ws.Cells.Item(24, 20).Value = 2.66
$ws.Cells.Item(24, 21).Value = '19/08/2023 16:59'
$ws.Cells.Item(24, 22).Value = 'https://www.betexplorer.com/football/poland/iii-liga-group-ii/swinoujscie-blekitni-stargard/SO98EN2g/'

# Row 25
$ws.Cells.Item(25, 6).Value = 'Elana Torun'
$ws.Cells.Item(25, 7).Value = 1
$ws.Cells.Item(25, 8).Value = 'Vineta W.'
$ws.Cells.Item(25, 9).Value = 2
$ws.Cells.Item(25, 10).Value = 1.86
$ws.Cells.Item(25, 11).Value = '19/08/2023 09:25'
$ws.Cells.Item(25, 12).Value = 1.83
$ws.Cells.Item(25, 13).Value = '19/08/2023 14:32'
$ws.Cells.Item(25, 14).Value = 3.54
$ws.Cells.Item(25, 15).Value = '19/08/2023 09:25'
$ws.Cells.Item(25, 16).Value = 3.49
$ws.Cells.Item(25, 17).Value = '19/08/2023 15:45'
$ws.Cells.Item(25, 18).Value = 3.41
$ws.Cells.Item(25, 19).Value = '19/08/2023 09:25'
$ws.Cells.Item(25, 20).Value = 3.67
$ws.Cells.Item(25, 21).Value = '19/08/2023 15:45'
$ws.Cells.Item(25, 22).Value = 'https://www.betexplorer.com/football/poland/iii-liga-group-ii/elana-torun-vineta-wolin/69LaGqYt/'

# Row 29
$ws.Cells.Item(29, 6).Value = 'Gedania Gdansk'
$ws.Cells.Item(29, 7).Value = 1
$ws.Cells.Item(29, 8).Value = 'Notec Czarnkow'
$ws.Cells.Item(29, 9).Value = 1
$ws.Cells.Item(29, 10).Value = 1.79
$ws.Cells.Item(29, 11).Value = '23/08/2023 11:12'
$ws.Cells.Item(29, 12).Value = 1.8
$ws.Cells.Item(29, 13).Value = '23/08/2023 17:28'
$ws.Cells.Item(29, 14).Value = 3.67
$ws.Cells.Item(29, 15).Value = '23/08/2023 11:12'
$ws.Cells.Item(29, 16).Value = 3.83
$ws.Cells.Item(29, 17).Value = '23/08/2023 17:28'
$ws.Cells.Item(29, 18).Value = 3.43
$ws.Cells.Item(29, 19).Value = '23/08/2023 11:12'
$ws.Cells.Item(29, 20).Value = 3.45
$ws.Cells.Item(29, 21).Value = '23/08/2023 17:28'
$ws.Cells.Item(29, 22).Value = 'https://www.betexplorer.com/football/poland/iii-liga-group-ii/gedania-gdansk-notec-czarnkow/zBliPvHB/'

# Row 30
$ws.Cells.Item(30, 6).Value = 'Kleczew'
$ws.Cells.Item(30, 7).Value = 4
$ws.Cells.Item(30, 8).Value = 'Starogard Gdanski'
$ws.Cells.Item(30, 9).Value = 3
$ws.Cells.Item(30, 10).Value = 2.12
$ws.Cells.Item(30, 11).Value = '22/08/2023 05:42'
$ws.Cells.Item(30, 12).Value = 1.98
$ws.Cells.Item(30, 13).Value = '23/08/2023 17:08'
$ws.Cells.Item(30, 14).Value = 3.25
$ws.Cells.Item(30, 15).Value = '22/08/2023 05:42'
$ws.Cells.Item(30, 16).Value = 4
$ws.Cells.Item(30, 17).Value = '23/08/2023 17:07'
$ws.Cells.Item(30, 18).Value = 2.72
$ws.Cells.Item(30, 19).Value = '22/08/2023 05:42'
$ws.Cells.Item(30, 20).Value = 2.85
$ws.Cells.Item(30, 21).Value = '23/08/2023 17:08'
$ws.Cells.Item(30, 22).Value = 'https://www.betexplorer.com/football/poland/iii-liga-group-ii/kleczew-starogard-gdanski/fiTHAu2I/'

# Row 32
$ws.Cells.Item(32, 6).Value = 'Vineta W.'
$ws.Cells.Item(32, 7).Value = 2
$ws.Cells.Item(32, 8).Value = 'Solec Kujawski'
$ws.Cells.Item(32, 9).Value = 1
$ws.Cells.Item(32, 10).Value = 1.76
$ws.Cells.Item(32, 11).Value = '22/08/2023 06:12'
$ws.Cells.Item(32, 12).Value = 1.74
$ws.Cells.Item(32, 13).Value = '23/08/2023 17:00'
$ws.Cells.Item(32, 14).Value = 3.54
$ws.Cells.Item(32, 15).Value = '22/08/2023 06:12'
$ws.Cells.Item(32, 16).Value = 3.85
$ws.Cells.Item(32, 17).Value = '23/08/2023 17:00'
$ws.Cells.Item(32, 18).Value = 3.34
$ws.Cells.Item(32, 19).Value = '22/08/2023 06:12'
$ws.Cells.Item(32, 20).Value = 3.68
$ws.Cells.Item(32, 21).Value = '23/08/2023 17:00'
$ws.Cells.Item(32, 22).Value = 'https://www.betexplorer.com/football/poland/iii-liga-group-ii/vineta-wolin-unia-solec-kujawski/pMevStXh/'

# Row 33
$ws.Cells.Item(33, 6).Value = 'Swit Skolwin'
$ws.Cells.Item(33, 7).Value = 3
$ws.Cells.Item(33, 8).Value = 'Luzino'
$ws.Cells.Item(33, 9).Value = 0
$ws.Cells.Item(33, 10).Value = 1.2
$ws.Cells.Item(33, 11).Value = '23/08/2023 11:12'
$ws.Cells.Item(33, 12).Value = 1.28
$ws.Cells.Item(33, 13).Value = '23/08/2023 17:50'
$ws.Cells.Item(33, 14).Value = 6.33
$ws.Cells.Item(33, 15).Value = '23/08/2023 11:12'
$ws.Cells.Item(33, 16).Value = 5.76
$ws.Cells.Item(33, 17).Value = '23/08/2023 17:50'
$ws.Cells.Item(33, 18).Value = 7.73
$ws.Cells.Item(33, 19).Value = '23/08/2023 11:12'
$ws.Cells.Item(33, 20).Value = 6.55
$ws.Cells.Item(33, 21).Value = '23/08/2023 17:50'
$ws.Cells.Item(33, 22).Value = 'https://www.betexplorer.com/football/poland/iii-liga-group-ii/swit-skolwin-luzino/hnanQK15/'

# Row 34
$ws.Cells.Item(34, 6).Value = 'Stolem Gniewino'
$ws.Cells.Item(34, 7).Value = 0
$ws.Cells.Item(34, 8).Value = 'Cartusia Kartuzy'
$ws.Cells.Item(34, 9).Value = 0
$ws.Cells.Item(34, 10).Value = 2.62
$ws.Cells.Item(34, 11).Value = '22/08/2023 06:12'
$ws.Cells.Item(34, 12).Value = 2.64
$ws.Cells.Item(34, 13).Value = '23/08/2023 17:59'
$ws.Cells.Item(34, 14).Value = 3.11
$ws.Cells.Item(34, 15).Value = '22/08/2023 06:12'
$ws.Cells.Item(34, 16).Value = 3.56
$ws.Cells.Item(34, 17).Value = '23/08/2023 17:59'
$ws.Cells.Item(34, 18).Value = 2.26
$ws.Cells.Item(34, 19).Value = '22/08/2023 06:12'
$ws.Cells.Item(34, 20).Value = 2.25
$ws.Cells.Item(34, 21).Value = '23/08/2023 17:59'
$ws.Cells.Item(34, 22).Value = 'https://www.betexplorer.com/football/poland/iii-liga-group-ii/stolem-gniewino-cartusia-kartuzy/noG3YrfP/'

# Row 35
$ws.Cells.Item(35, 6).Value = 'Nowe Skalmierzyce'
$ws.Cells.Item(35, 7).Value = 0
$ws.Cells.Item(35, 8).Value = 'Pogon Szczecin II'
$ws.Cells.Item(35, 9).Value = 5
$ws.Cells.Item(35, 10).Value = 2.6
$ws.Cells.Item(35, 11).Value = '22/08/2023 06:12'
$ws.Cells.Item(35, 12).Value = 1.95
$ws.Cells.Item(35, 13).Value = '23/08/2023 17:45'
$ws.Cells.Item(35, 14).Value = 3.2
$ws.Cells.Item(35, 15).Value = '22/08/2023 06:12'
$ws.Cells.Item(35, 16).Value = 3.4
$ws.Cells.Item(35, 17).Value = '23/08/2023 17:45'
$ws.Cells.Item(35, 18).Value = 2.22
$ws.Cells.Item(35, 19).Value = '22/08/2023 06:12'
$ws.Cells.Item(35, 20).Value = 2.9
$ws.Cells.Item(35, 21).Value = '23/08/2023 17:45'
$ws.Cells.Item(35, 22).Value = 'https://www.betexplorer.com/football/poland/iii-liga-group-ii/nowe-skalmierzyce-pogon-szczecin/YsSL9aHO/'

# Row 40
$ws.Cells.Item(40, 6).Value = 'Luzino'
$ws.Cells.Item(40, 7).Value = 1
$ws.Cells.Item(40, 8).Value = 'Gedania Gdansk'
$ws.Cells.Item(40, 9).Value = 2
$ws.Cells.Item(40, 10).Value = 4.25
$ws.Cells.Item(40, 11).Value = '27/08/2023 03:42'
$ws.Cells.Item(40, 12).Value = 3
$ws.Cells.Item(40, 13).Value = '27/08/2023 10:03'
$ws.Cells.Item(40, 14).Value = 4.61
$ws.Cells.Item(40, 15).Value = '27/08/2023 03:42'
$ws.Cells.Item(40, 16).Value = 4.47
$ws.Cells.Item(40, 17).Value = '27/08/2023 10:08'
$ws.Cells.Item(40, 18).Value = 1.51
$ws.Cells.Item(40, 19).Value = '27/08/2023 03:42'
$ws.Cells.Item(40, 20).Value = 1.82
$ws.Cells.Item(40, 21).Value = '27/08/2023 10:03'
$ws.Cells.Item(40, 22).Value = 'https://www.betexplorer.com/football/poland/iii-liga-group-ii/luzino-gedania-gdansk/SI6bNIoO/'

# Row 41
$ws.Cells.Item(41, 6).Value = 'Pogon Szczecin II'
$ws.Cells.Item(41, 7).Value = 2
$ws.Cells.Item(41, 8).Value = 'Stolem Gniewino'
$ws.Cells.Item(41, 9).Value = 0
$ws.Cells.Item(41, 10).Value = 1.35
$ws.Cells.Item(41, 11).Value = '25/08/2023 23:13'
$ws.Cells.Item(41, 12).Value = 1.41
$ws.Cells.Item(41, 13).Value = '27/08/2023 10:56'
$ws.Cells.Item(41, 14).Value = 4.37
$ws.Cells.Item(41, 15).Value = '25/08/2023 23:13'
$ws.Cells.Item(41, 16).Value = 4.81
$ws.Cells.Item(41, 17).Value = '27/08/2023 10:57'
$ws.Cells.Item(41, 18).Value = 5.45
$ws.Cells.Item(41, 19).Value = '25/08/2023 23:13'
$ws.Cells.Item(41, 20).Value = 5.27
$ws.Cells.Item(41, 21).Value = '27/08/2023 10:56'
$ws.Cells.Item(41, 22).Value = 'https://www.betexplorer.com/football/poland/iii-liga-group-ii/pogon-szczecin-stolem-gniewino/zPTPwvW4/'

# Row 46
$ws.Cells.Item(46, 6).Value = 'Swit Skolwin'
$ws.Cells.Item(46, 7).Value = 3
$ws.Cells.Item(46, 8).Value = 'Solec Kujawski'
$ws.Cells.Item(46, 9).Value = 1
$ws.Cells.Item(46, 10).Value = 1.36
$ws.Cells.Item(46, 11).Value = '31/08/2023 05:42'
$ws.Cells.Item(46, 12).Value = 1.31
$ws.Cells.Item(46, 13).Value = '01/09/2023 16:39'
$ws.Cells.Item(46, 14).Value = 4.45
$ws.Cells.Item(46, 15).Value = '31/08/2023 05:42'
$ws.Cells.Item(46, 16).Value = 5.15
$ws.Cells.Item(46, 17).Value = '01/09/2023 17:00'
$ws.Cells.Item(46, 18).Value = 5.19
$ws.Cells.Item(46, 19).Value = '31/08/2023 05:42'
$ws.Cells.Item(46, 20).Value = 6.86
$ws.Cells.Item(46, 21).Value = '01/09/2023 17:00'
$ws.Cells.Item(46, 22).Value = 'https://www.betexplorer.com/football/poland/iii-liga-group-ii/swit-skolwin-unia-solec-kujawski/Uokahwob/'

# Row 47
$ws.Cells.Item(47, 6).Value = 'Vineta W.'
$ws.Cells.Item(47, 7).Value = 2
$ws.Cells.Item(47, 8).Value = 'Unia Swarzedz'
$ws.Cells.Item(47, 9).Value = 1
$ws.Cells.Item(47, 10).Value = 2.21
$ws.Cells.Item(47, 11).Value = '31/08/2023 05:42'
$ws.Cells.Item(47, 12).Value = 2.28
$ws.Cells.Item(47, 13).Value = '01/09/2023 17:00'
$ws.Cells.Item(47, 14).Value = 3.24
$ws.Cells.Item(47, 15).Value = '31/08/2023 05:42'
$ws.Cells.Item(47, 16).Value = 3.38
$ws.Cells.Item(47, 17).Value = '01/09/2023 17:00'
$ws.Cells.Item(47, 18).Value = 2.59
$ws.Cells.Item(47, 19).Value = '31/08/2023 05:42'
$ws.Cells.Item(47, 20).Value = 2.7
$ws.Cells.Item(47, 21).Value = '01/09/2023 17:00'
$ws.Cells.Item(47, 22).Value = 'https://www.betexplorer.com/football/poland/iii-liga-group-ii/vineta-wolin-unia-swarzedz/YJejfaGo/'

# Row 52
$ws.Cells.Item(52, 6).Value = 'Unia Swarzedz'
$ws.Cells.Item(52, 7).Value = 1
$ws.Cells.Item(52, 8).Value = 'Sroda'
$ws.Cells.Item(52, 9).Value = 1
$ws.Cells.Item(52, 10).Value = 2.35
$ws.Cells.Item(52, 11).Value = '07/09/2023 05:12'
$ws.Cells.Item(52, 12).Value = 1.77
$ws.Cells.Item(52, 13).Value = '08/09/2023 16:56'
$ws.Cells.Item(52, 14).Value = 3.27
$ws.Cells.Item(52, 15).Value = '07/09/2023 05:12'
$ws.Cells.Item(52, 16).Value = 3.76
$ws.Cells.Item(52, 17).Value = '08/09/2023 16:56'
$ws.Cells.Item(52, 18).Value = 2.41
$ws.Cells.Item(52, 19).Value = '07/09/2023 05:12'
$ws.Cells.Item(52, 20).Value = 3.64
$ws.Cells.Item(52, 21).Value = '08/09/2023 16:56'
$ws.Cells.Item(52, 22).Value = 'https://www.betexplorer.com/football/poland/iii-liga-group-ii/unia-swarzedz-polonia-sroda-wielkopol/tl7iV6q2/'

# Row 53
$ws.Cells.Item(53, 6).Value = 'Luzino'
$ws.Cells.Item(53, 7).Value = 0
$ws.Cells.Item(53, 8).Value = 'Kleczew'
$ws.Cells.Item(53, 9).Value = 0
$ws.Cells.Item(53, 10).Value = 3.31
$ws.Cells.Item(53, 11).Value = '08/09/2023 07:42'
$ws.Cells.Item(53, 12).Value = 3.42
$ws.Cells.Item(53, 13).Value = '08/09/2023 16:57'
$ws.Cells.Item(53, 14).Value = 3.77
$ws.Cells.Item(53, 15).Value = '08/09/2023 07:42'
$ws.Cells.Item(53, 16).Value = 4.02
$ws.Cells.Item(53, 17).Value = '08/09/2023 16:58'
$ws.Cells.Item(53, 18).Value = 1.83
$ws.Cells.Item(53, 19).Value = '08/09/2023 07:42'
$ws.Cells.Item(53, 20).Value = 1.73
$ws.Cells.Item(53, 21).Value = '08/09/2023 16:58'
$ws.Cells.Item(53, 22).Value = 'https://www.betexplorer.com/football/poland/iii-liga-group-ii/luzino-kleczew/EZjBkyVG/'

# --- Part 2: append two new match rows (70 and 71) ---

# Row 70
$ws.Cells.Item(69, 1).Copy() | Out-Null
$ws.Cells.Item(70, 1).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(69, 5).Copy() | Out-Null
$ws.Cells.Item(70, 5).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(70, 1).Value = 69
$ws.Cells.Item(70, 2).Value = 'poland'
$ws.Cells.Item(70, 3).Value = 'iii-liga-group-ii'
$ws.Cells.Item(70, 4).Value = '2023-2024'
$ws.Cells.Item(70, 5).Value = 45191.5
$ws.Cells.Item(70, 6).Value = 'Pogon Szczecin II'
$ws.Cells.Item(70, 7).Value = 3
$ws.Cells.Item(70, 8).Value = 'Sroda'
$ws.Cells.Item(70, 9).Value = 0
$ws.Cells.Item(70, 10).Value = 1.59
$ws.Cells.Item(70, 11).Value = '20/09/2023 23:13'
$ws.Cells.Item(70, 12).Value = 1.45
$ws.Cells.Item(70, 13).Value = '22/09/2023 11:58'
$ws.Cells.Item(70, 14).Value = 3.86
$ws.Cells.Item(70, 15).Value = '20/09/2023 23:13'
$ws.Cells.Item(70, 16).Value = 4.53
$ws.Cells.Item(70, 17).Value = '22/09/2023 11:58'
$ws.Cells.Item(70, 18).Value = 3.82
$ws.Cells.Item(70, 19).Value = '20/09/2023 23:13'
$ws.Cells.Item(70, 20).Value = 5.14
$ws.Cells.Item(70, 21).Value = '22/09/2023 11:58'
$ws.Cells.Item(70, 22).Value = 'https://www.betexplorer.com/football/poland/iii-liga-group-ii/pogon-szczecin-polonia-sroda-wielkopol/zoN4e3c7/'

# Row 71
$ws.Cells.Item(69, 1).Copy() | Out-Null
$ws.Cells.Item(71, 1).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(69, 5).Copy() | Out-Null
$ws.Cells.Item(71, 5).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(71, 1).Value = 70
$ws.Cells.Item(71, 2).Value = 'poland'
$ws.Cells.Item(71, 3).Value = 'iii-liga-group-ii'
$ws.Cells.Item(71, 4).Value = '2023-2024'
$ws.Cells.Item(71, 5).Value = 45191.70833333334
$ws.Cells.Item(71, 6).Value = 'Unia Swarzedz'
$ws.Cells.Item(71, 7).Value = 3
$ws.Cells.Item(71, 8).Value = 'Gedania Gdansk'
$ws.Cells.Item(71, 9).Value = 4
$ws.Cells.Item(71, 10).Value = 2.11
$ws.Cells.Item(71, 11).Value = '21/09/2023 04:12'
$ws.Cells.Item(71, 12).Value = 2.19
$ws.Cells.Item(71, 13).Value = '22/09/2023 16:59'
$ws.Cells.Item(71, 14).Value = 3.36
$ws.Cells.Item(71, 15).Value = '21/09/2023 04:12'
$ws.Cells.Item(71, 16).Value = 3.3
$ws.Cells.Item(71, 17).Value = '22/09/2023 16:59'
$ws.Cells.Item(71, 18).Value = 2.65
$ws.Cells.Item(71, 19).Value = '21/09/2023 04:12'
$ws.Cells.Item(71, 20).Value = 2.9
$ws.Cells.Item(71, 21).Value = '22/09/2023 16:59'
$ws.Cells.Item(71, 22).Value = 'https://www.betexplorer.com/football/poland/iii-liga-group-ii/unia-swarzedz-gedania-gdansk/2XOdcPSf/'

$excel.CutCopyMode = 0
